$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename "Activity" column to "Description" ---
$ws.Range("B1").Value = "Description"

# --- Replace the 15 "Activity" entries with new "Description" entries ---
$descriptions = @(
    "Initial concept sketches for new attraction area ",
    "Site layout planning and space allocation ",
    "Environmental impact assessment of new structures ",
    "Character meet-and-greet zone design ",
    "Traffic flow analysis for expanded visitor capacity ",
    "Facade detailing for themed restaurant ",
    "Water feature integration planning ",
    "Lighting design for nighttime operations ",
    "Visitor experience impact study for construction period ",
    "Foundation assessment for main attraction ",
    "Guest flow optimization for new pathways ",
    "Theming coordination with creative team ",
    "Accessibility compliance review ",
    "Safety protocol documentation ",
    "Final design review with client "
)

for ($i = 0; $i -lt $descriptions.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $descriptions[$i]
}

# --- Correct a handful of the timesheet dates in column A ---
$ws.Range("A3").Value = 45292
$ws.Range("A5").Value = 45294
$ws.Range("A10").Value = 45299
$ws.Range("A13").Value = 45302
$ws.Range("A14").Value = 45302

# --- Update the active selection ---
$ws.Range("I18").Select() | Out-Null

# --- Force portrait page orientation (adds a pageSetup element) ---
$ws.PageSetup.Orientation = 1
